$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Paragraph "2024.8.28   天气晴": merge the two runs (which used to
#    be split by the _GoBack bookmark) into a single run, and drop the
#    bookmark from this location (it moves further down, see step 3).
# ------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$dateParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.StartsWith("2024.8.28")) {
        $dateParaIndex = $i
        break
    }
}

$dateRange = $d.Range($d.Paragraphs.Item($dateParaIndex).Range.Start, $d.Paragraphs.Item($dateParaIndex).Range.End - 1)
# Re-assert the same visible text (with one extra sentinel char) so the
# COM layer actually rewrites the range as a single run instead of
# treating an identical assignment as a no-op; then trim the sentinel.
$dateRange.Text = "2024.8.28   天气晴#"
$dateParaAgain = $d.Paragraphs.Item($dateParaIndex)
$sentinel = $d.Range($dateParaAgain.Range.End - 2, $dateParaAgain.Range.End - 1)
$sentinel.Delete()

# ------------------------------------------------------------------
# 2) Diary paragraph: fix the "bind;tap" typo to "bind:tap" and expand
#    the closing sentence.
# ------------------------------------------------------------------
$ok1 = $d.Content.Find.Execute("bind;tap", $true, $false, $false, $false, $false, $true, 1, $false, "bind:tap", 2)
Write-Output "fix bind;tap -> bind:tap : $ok1"

$ok2 = $d.Content.Find.Execute("我计划继续进行菜谱页面的完善，使其更加美观。", $true, $false, $false, $false, $false, $true, 1, $false, "我计划继续进行菜谱页面的美观优化以及排版设计，使其更加美观耐看、引人注目。", 2)
Write-Output "expand closing sentence : $ok2"

# ------------------------------------------------------------------
# 3) Split the diary run right after "bind:" and drop the _GoBack
#    bookmark back in at that split point.
# ------------------------------------------------------------------
$findRange = $d.Content
$ok3 = $findRange.Find.Execute("bind:tap")
Write-Output "locate bind:tap split point : $ok3"
$splitPos = $findRange.Start + 5

$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
